$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (font/border/alignment) from an existing header cell so the
# new headers look consistent with the rest of row 1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the header text since PasteSpecial(formats) shouldn't touch values,
# but set again to be safe
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the new data row values
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
